$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.665929
$ws.Range("H2").Value = 46.997787
$ws.Range("I2").Value = 0.5304105216372652
$ws.Range("J2").Value = 0.5304105216372653
$ws.Range("M2").Value = 37.641477
$ws.Range("N2").Value = 112.924431
$ws.Range("O2").Value = 0.4464054516632745
$ws.Range("P2").Value = 0.4464054516632746
$ws.Range("Q2").Value = 589.6887061371331
$ws.Range("R2").Value = 5307.198355234197
$ws.Range("S2").Value = 0.2367781484784364
$ws.Range("T2").Value = 0.2367781484784365

$ws.Range("G3").Value = 15.665929
$ws.Range("H3").Value = 46.997787
$ws.Range("I3").Value = 0.5304105216372652
$ws.Range("J3").Value = 0.5304105216372653
$ws.Range("N3").Value = 0.996564
$ws.Range("O3").Value = 0.003939551420288843
$ws.Range("P3").Value = 0.003939551420288844
$ws.Range("Q3").Value = 5.204033622652
$ws.Range("R3").Value = 46.83630260386801
$ws.Range("S3").Value = 0.002089579523852234
$ws.Range("T3").Value = 0.002089579523852235

$ws.Range("G4").Value = 15.665929
$ws.Range("H4").Value = 46.997787
$ws.Range("I4").Value = 0.5304105216372652
$ws.Range("J4").Value = 0.5304105216372653
$ws.Range("M4").Value = 46.34761033333334
$ws.Range("N4").Value = 139.042831
$ws.Range("O4").Value = 0.5496549969164365
$ws.Range("P4").Value = 0.5496549969164366
$ws.Range("Q4").Value = 726.0783728016664
$ws.Range("R4").Value = 6534.705355214997
$ws.Range("S4").Value = 0.2915427936349765
$ws.Range("T4").Value = 0.2915427936349766

$ws.Range("I5").Value = 0.01939778913130443
$ws.Range("J5").Value = 0.01939778913130443
$ws.Range("M5").Value = 37.641477
$ws.Range("N5").Value = 112.924431
$ws.Range("O5").Value = 0.4464054516632745
$ws.Range("P5").Value = 0.4464054516632746
$ws.Range("Q5").Value = 21.565667927271
$ws.Range("R5").Value = 194.091011345439
$ws.Range("S5").Value = 0.008659278818428911
$ws.Range("T5").Value = 0.008659278818428913

$ws.Range("I6").Value = 0.01939778913130443
$ws.Range("J6").Value = 0.01939778913130443
$ws.Range("N6").Value = 0.996564
$ws.Range("O6").Value = 0.003939551420288843
$ws.Range("P6").Value = 0.003939551420288844
$ws.Range("S6").Value = 0.00007641858772269386
$ws.Range("T6").Value = 0.00007641858772269388

$ws.Range("I7").Value = 0.01939778913130443
$ws.Range("J7").Value = 0.01939778913130443
$ws.Range("M7").Value = 46.34761033333334
$ws.Range("N7").Value = 139.042831
$ws.Range("O7").Value = 0.5496549969164365
$ws.Range("P7").Value = 0.5496549969164366
$ws.Range("Q7").Value = 26.55361195500433
$ws.Range("R7").Value = 238.982507595039
$ws.Range("S7").Value = 0.01066209172515282
$ws.Range("T7").Value = 0.01066209172515283

$ws.Range("G8").Value = 13.29662733333333
$ws.Range("H8").Value = 39.889882
$ws.Range("I8").Value = 0.4501916892314303
$ws.Range("J8").Value = 0.4501916892314304
$ws.Range("M8").Value = 37.641477
$ws.Range("N8").Value = 112.924431
$ws.Range("O8").Value = 0.4464054516632745
$ws.Range("P8").Value = 0.4464054516632746
$ws.Range("Q8").Value = 500.5046919452381
$ws.Range("R8").Value = 4504.542227507142
$ws.Range("S8").Value = 0.2009680243664092
$ws.Range("T8").Value = 0.2009680243664092

$ws.Range("G9").Value = 13.29662733333333
$ws.Range("H9").Value = 39.889882
$ws.Range("I9").Value = 0.4501916892314303
$ws.Range("J9").Value = 0.4501916892314304
$ws.Range("N9").Value = 0.996564
$ws.Range("O9").Value = 0.003939551420288843
$ws.Range("P9").Value = 0.003939551420288844
$ws.Range("Q9").Value = 4.416980040605333
$ws.Range("R9").Value = 39.752820365448
$ws.Range("S9").Value = 0.001773553308713915
$ws.Range("T9").Value = 0.001773553308713916

$ws.Range("G10").Value = 13.29662733333333
$ws.Range("H10").Value = 39.889882
$ws.Range("I10").Value = 0.4501916892314303
$ws.Range("J10").Value = 0.4501916892314304
$ws.Range("M10").Value = 46.34761033333334
$ws.Range("N10").Value = 139.042831
$ws.Range("O10").Value = 0.5496549969164365
$ws.Range("P10").Value = 0.5496549969164366
$ws.Range("Q10").Value = 616.2669023928825
$ws.Range("R10").Value = 5546.402121535943
$ws.Range("S10").Value = 0.2474501115563072
$ws.Range("T10").Value = 0.2474501115563073

